$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists one row per state/territory (NAME in column B). The commit
# removes the "District of Columbia" and "Puerto Rico" rows entirely (and,
# as a consequence, their "DC"/"PR" abbreviation strings) from this
# work-economic relative-score table, shifting every following row up.
#
# Locate each row by its NAME value instead of a hard-coded row number so
# the deletions are resilient to row position, then delete the entire row
# (Excel shifts the rows below upward and updates the used range / shared
# strings automatically).

$dcCell = $ws.Columns.Item(2).Find("District of Columbia")
if ($dcCell -ne $null) {
    $ws.Rows.Item($dcCell.Row).Delete()
}

$prCell = $ws.Columns.Item(2).Find("Puerto Rico")
if ($prCell -ne $null) {
    $ws.Rows.Item($prCell.Row).Delete()
}
